# spec/fixtures/invalid_customer_import.xlsx
# Rename the three "Zip Code" column headers to "Zipcode" (Billing, Primary
# Shipping, Secondary Shipping) and update the saved selection/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Billing Zipcode"
$ws.Range("J1").Value = "Primary Shipping Zipcode"
$ws.Range("N1").Value = "Secondary Shipping Zipcode"

# Match the committed cursor/selection position.
$ws.Range("H15").Select() | Out-Null
